$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Range("E2").Value = 17
$ws.Range("E6").Value = 39
$ws.Range("E7").Value = 17
$ws.Range("E12").Value = 17
$ws.Range("E15").Value = 65
$ws.Range("E16").Value = 230
$ws.Range("E18").Value = 63
